# ---------------------------------------------------------------------------
# Adds the "COTRIFRED" cooperative's results to the three reporting sheets:
#   - comparativo_master : 16 new detail rows (352-367)
#   - niveis_master      : 3 new detail rows (56-58) ahead of the TOTAL rows,
#                           which shift down and get their totals updated
#   - financeiro_master  : a new detail row (COTRIFRED) replacing the old
#                           TOTAL row, followed by a refreshed TOTAL row
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1) comparativo_master (sheet1) - append 16 rows at the bottom
# ===========================================================================
$wsComp = $wb.Worksheets.Item("comparativo_master")

$compRows = @(
    @('COTRIFRED','ALCIONE MACIAK',11,'Básico',47,'Intermediário',36,327.27272727272731),
    @('COTRIFRED','ANDERSON LUIS TUR',16,'Básico',80,'Avançado',64,400),
    @('COTRIFRED','DANIEL PINHEIRO',14,'Básico',80,'Avançado',66,471.42857142857139),
    @('COTRIFRED','FERNANDO FELIN',21,'Básico',61,'Avançado',40,190.47619047619051),
    @('COTRIFRED','IONARA FATIMA FRIZON PIOVESAN',16,'Básico',80,'Avançado',64,400),
    @('COTRIFRED','JOAO LOPES',19,'Básico',81,'Avançado',62,326.31578947368422),
    @('COTRIFRED','JOARES MAGALSKI',4,'Básico',81,'Avançado',77,1925),
    @('COTRIFRED','LUCAS BELMIRO AZEVEDO',4,'Básico',60,'Avançado',56,1400),
    @('COTRIFRED','LUCAS RISSI PIOVESAN',11,'Básico',72,'Avançado',61,554.54545454545462),
    @('COTRIFRED','MATHEUS EUGENIO SARI',6,'Básico',83,'Avançado',77,1283.333333333333),
    @('COTRIFRED','ROSELENE SARI BEVILACOUA',18,'Básico',86,'Avançado',68,377.77777777777783),
    @('COTRIFRED','SANDRO ROBERTO CAPOANE',2,'Básico',44,'Intermediário',42,2100),
    @('COTRIFRED','TIAGO STIVAL',19,'Básico',54,'Intermediário',35,184.21052631578951),
    @('COTRIFRED','VALDINEI MAZZONETTO',14,'Básico',64,'Avançado',50,357.14285714285722),
    @('COTRIFRED','VANIO MIOR',19,'Básico',79,'Avançado',60,315.78947368421052),
    @('COTRIFRED','VOLNEI ZONTA',8,'Básico',86,'Avançado',78,975)
)

$startRow = 352
for ($i = 0; $i -lt $compRows.Count; $i++) {
    $r = $startRow + $i
    $row = $compRows[$i]
    $wsComp.Cells.Item($r,1).Value = $row[0]
    $wsComp.Cells.Item($r,2).Value = $row[1]
    $wsComp.Cells.Item($r,3).Value = $row[2]
    $wsComp.Cells.Item($r,4).Value = $row[3]
    $wsComp.Cells.Item($r,5).Value = $row[4]
    $wsComp.Cells.Item($r,6).Value = $row[5]
    $wsComp.Cells.Item($r,7).Value = $row[6]
    $wsComp.Cells.Item($r,8).Value = $row[7]
}

# ===========================================================================
# 2) niveis_master (sheet3) - insert 3 rows before the TOTAL rows and
#    refresh the TOTAL values
# ===========================================================================
$wsNiv = $wb.Worksheets.Item("niveis_master")

$wsNiv.Rows.Item(56).Insert()
$wsNiv.Rows.Item(56).Insert()
$wsNiv.Rows.Item(56).Insert()

$wsNiv.Cells.Item(56,1).Value = "COTRIFRED"
$wsNiv.Cells.Item(56,2).Value = "Básico"
$wsNiv.Cells.Item(56,3).Value = 16
$wsNiv.Cells.Item(56,4).Value = 0

$wsNiv.Cells.Item(57,1).Value = "COTRIFRED"
$wsNiv.Cells.Item(57,2).Value = "Intermediário"
$wsNiv.Cells.Item(57,3).Value = 0
$wsNiv.Cells.Item(57,4).Value = 3

$wsNiv.Cells.Item(58,1).Value = "COTRIFRED"
$wsNiv.Cells.Item(58,2).Value = "Avançado"
$wsNiv.Cells.Item(58,3).Value = 0
$wsNiv.Cells.Item(58,4).Value = 13

# Rows 59-61 are the old TOTAL rows (56-58), shifted down by the insert;
# refresh their counts to include the new COTRIFRED entries.
$wsNiv.Cells.Item(59,3).Value = 3
$wsNiv.Cells.Item(59,4).Value = 210

$wsNiv.Cells.Item(60,3).Value = 301
$wsNiv.Cells.Item(60,4).Value = 34

$wsNiv.Cells.Item(61,3).Value = 62
$wsNiv.Cells.Item(61,4).Value = 122

# Refresh the AutoFilter range (A1:D58 -> A1:D61) without leaving the filter
# toggled off, and update the matching hidden defined name.
$wsNiv.AutoFilterMode = $false
$wsNiv.Range("A1:D61").AutoFilter()

$nivName = $wb.Names.Item("niveis_master!_FilterDatabase")
$nivName.RefersTo = "=niveis_master!`$A`$1:`$D`$61"

# ===========================================================================
# 3) financeiro_master (sheet4) - turn old TOTAL row into the COTRIFRED row
#    and append a fresh TOTAL row
# ===========================================================================
$wsFin = $wb.Worksheets.Item("financeiro_master")

$wsFin.Cells.Item(20,1).Value = "COTRIFRED"
$wsFin.Cells.Item(20,2).Value = "Gestão Financeira"
$wsFin.Cells.Item(20,3).Value = 7
$wsFin.Cells.Item(20,4).Value = 171
$wsFin.Cells.Item(20,5).Value = 164
$wsFin.Cells.Item(20,6).Value = 2342.8571428571431

$wsFin.Cells.Item(21,1).Value = "TOTAL"
$wsFin.Cells.Item(21,2).Value = "Gestão Financeira"
$wsFin.Cells.Item(21,3).Value = 517
$wsFin.Cells.Item(21,4).Value = 3667
$wsFin.Cells.Item(21,5).Value = 3134
$wsFin.Cells.Item(21,6).Value = 516.92913385826773

# ===========================================================================
# 4) Restore / update on-screen selections for each sheet
# ===========================================================================
$wsNiv.Activate()
$wsNiv.Range("B1").Select()

$wsFin.Activate()
$wsFin.Range("F21").Select()

$wsComp.Activate()
$wsComp.Range("A353:A367").Select()
